$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "posibles respuestas" header (column D) to "posibles_respuestas"
$ws.Range("D1").Value = "posibles_respuestas"

# Reflect the resulting view state: selection moved to D2, sheet scrolled right
# so column C becomes the left-most visible column.
$ws.Range("D2").Select()
$excel.ActiveWindow.ScrollColumn = 3
